# Update column F (dSF) values for several rows to repull data / push all data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -7
    4  = -7
    5  = -2
    6  = -7
    11 = 1
    13 = 10
    14 = -8
    15 = 0
    16 = 0
    18 = -6
    23 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
